$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item(3)   # "Sheet3"
$ws4 = $wb.Worksheets.Item(4)   # "Sheet4"

# --- Data edits on Sheet3 -------------------------------------------------
# Rx-1 Tx measurement (B3): 177 -> 177.5
$ws3.Range("B3").Value = 177.5
# Rx-2 Tx measurement (B4): 178 -> 177.5
$ws3.Range("B4").Value = 177.5

# --- Sheet view / selection updates --------------------------------------
# Sheet4 previously had the tab selected and A15:D18 selected; move the
# selection back to a single cell (E16) and drop the tab-selected flag by
# activating Sheet3 instead.
[void]$ws4.Range("E16").Select()

# Sheet3 becomes the active tab, with A1:D18 selected (active cell ends on
# the last cell of the drag, D18).
[void]$ws3.Range("A1:D18").Select()
[void]$ws3.Activate()
